# Updates cryptos list: refresh Price and Volume(1h) columns for most rows,
# and for rows 44/45 swap EnergySwap/VeChain order with new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.673.46"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.474.94"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "318.85"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "93.07"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +8.57%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "2.855.76"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("E14").Value = "  +0.94%  "
$ws.Range("D15").Value = "15.76"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").Value = "2.470.27"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "41.641.34"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "11.32"
$ws.Range("D23").Value = "239.63"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("E28").Value = "  +2.69%  "
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").Value = "36.17"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").Value = "158.53"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "2.58"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").Value = "17.34"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("D37").Value = "1.89"
$ws.Range("E37").Value = "  +5.57%  "
$ws.Range("D38").Value = "2.94"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  +5.10%  "
$ws.Range("D43").Value = "1.995.11"
$ws.Range("E43").Value = "  +2.17%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.04"
$ws.Range("E44").Value = "  +2.10%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("E47").Value = "  +3.00%  "
$ws.Range("D48").Value = "2.712.72"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "97.35"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "74.27"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("D51").Value = "67.15"
$ws.Range("E51").Value = "  +0.47%  "
